$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The LR-pair table (Fgf22 -> Fgfr2) now also includes "ECs" as a sending cluster
# (in addition to the existing "MuSCs"), so the 3 data rows become 6 data rows,
# each column A-T holding the updated TPM-derived values.
$rows = @(
 ,@("ECs","Fgf22","Fgfr2","ECs",1,0.3333333333333333,0.09830299999999999,0.294909,0.3722197960868258,0.3722197960868259,1,0.3333333333333333,0.036942,0.110826,0.02099032928903418,0.02099032928903418,0.003631509425999999,0.032683584834,0.007813016087759628,0.007813016087759632)
 ,@("ECs","Fgf22","Fgfr2","FAPs",1,0.3333333333333333,0.09830299999999999,0.294909,0.3722197960868258,0.3722197960868259,3,1,0.9431116666666667,2.829335,0.5358731102718634,0.5358731102718634,0.09271070616833332,0.8343963555149999,0.1994625798338061,0.1994625798338062)
 ,@("ECs","Fgf22","Fgfr2","MuSCs",1,0.3333333333333333,0.09830299999999999,0.294909,0.3722197960868258,0.3722197960868259,3,1,0.7798996666666667,2.339699,0.4431365604391025,0.4431365604391026,0.07666647693233332,0.689998292391,0.1649442001652601,0.1649442001652602)
 ,@("MuSCs","Fgf22","Fgfr2","ECs",2,0.6666666666666666,0.1657963333333334,0.497389,0.6277802039131741,0.6277802039131741,1,0.3333333333333333,0.036942,0.110826,0.02099032928903418,0.02099032928903418,0.006124848146,0.055123633314,0.01317731320127455,0.01317731320127455)
 ,@("MuSCs","Fgf22","Fgfr2","FAPs",2,0.6666666666666666,0.1657963333333334,0.497389,0.6277802039131741,0.6277802039131741,3,1,0.9431116666666667,2.829335,0.5358731102718634,0.5358731102718634,0.1563644562572222,1.407280106315,0.3364105304380572,0.3364105304380572)
 ,@("MuSCs","Fgf22","Fgfr2","MuSCs",2,0.6666666666666666,0.1657963333333334,0.497389,0.6277802039131741,0.6277802039131741,3,1,0.7798996666666667,2.339699,0.4431365604391025,0.4431365604391026,0.1293045051012222,1.163740545911,0.2781923602738424,0.2781923602738424)
)

$r = 2
foreach ($row in $rows) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c++
  }
  $r++
}
